# Update recalculated TPM-derived statistics in the LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 0.6142318033764959
$ws.Range("J2").Value = 0.6142318033764957
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 4.780548468954444
$ws.Range("R2").Value = 43.02493622059
$ws.Range("S2").Value = 0.02847271815842944
$ws.Range("T2").Value = 0.02847271815842944

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.6142318033764959
$ws.Range("J3").Value = 0.6142318033764957
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("S3").Value = 0.4245884951329531
$ws.Range("T3").Value = 0.424588495132953

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.6142318033764959
$ws.Range("J4").Value = 0.6142318033764957
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 27.06042371454333
$ws.Range("R4").Value = 243.54381343089
$ws.Range("S4").Value = 0.1611705900851133
$ws.Range("T4").Value = 0.1611705900851133

# Row 5 (MuSCs -> ECs)
$ws.Range("G5").Value = 1.102210333333334
$ws.Range("H5").Value = 3.306631
$ws.Range("I5").Value = 0.3857681966235041
$ws.Range("J5").Value = 0.3857681966235041
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 3.002422785017222
$ws.Range("R5").Value = 27.021805065155
$ws.Range("S5").Value = 0.01788228658393649
$ws.Range("T5").Value = 0.01788228658393649

# Row 6 (MuSCs -> FAPs)
$ws.Range("G6").Value = 1.102210333333334
$ws.Range("H6").Value = 3.306631
$ws.Range("I6").Value = 0.3857681966235041
$ws.Range("J6").Value = 0.3857681966235041
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("Q6").Value = 44.7724788673169
$ws.Range("R6").Value = 402.9523098058521
$ws.Range("S6").Value = 0.2666627438926821
$ws.Range("T6").Value = 0.266662743892682

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 1.102210333333334
$ws.Range("H7").Value = 3.306631
$ws.Range("I7").Value = 0.3857681966235041
$ws.Range("J7").Value = 0.3857681966235041
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 16.99529525961167
$ws.Range("R7").Value = 152.957657336505
$ws.Range("S7").Value = 0.1012231661468855
$ws.Range("T7").Value = 0.1012231661468855
